$d = $word.ActiveDocument

# 1. "(Ace: 1, King: 13). These VALUE, along with its SUIT, i..."
#    -> "(Ace: 1, King: 13). The VALUE, along with its SUIT, i..."
$d.Content.Find.Execute(
    "These VALUE, along with its SUIT",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The VALUE, along with its SUIT", 2) | Out-Null

# 2. Add a new trailing sentence about how/when effects can be activated,
#    right after "...they can be activated either turn."
$d.Content.Find.Execute(
    "Depending on the effects themselves, they can be activated either turn.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Depending on the effects themselves, they can be activated either turn. Effects can be activated in two ways: while the card is FACE-UP on the field, or DISCARDED to the graveyard.",
    2) | Out-Null

# 3. "Noone can attack in the first turn." -> "No one can attack in the first turn."
$d.Content.Find.Execute(
    "Noone can attack in the first turn.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "No one can attack in the first turn.", 2) | Out-Null

# 4. Collapse the extra spaces in "...targeted as an attack    target." down to a single space.
$d.Content.Find.Execute(
    "targeted as an attack    target.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "targeted as an attack target.", 2) | Out-Null

# 5. Move the (rendering-only) page-break hint from just before "J - K: " to
#    just before "below:" in the previous sentence. Touching the "J - K: " run
#    drops the stale marker there; re-touching the "...listed below:" run
#    establishes the new split point.
$d.Content.Find.Execute(
    "Cards with power <= 6 can be simply summoned. Higher ones will need sacrifices, listed below:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Cards with power <= 6 can be simply summoned. Higher ones will need sacrifices, listed below:",
    2) | Out-Null

$d.Content.Find.Execute(
    "J - K: ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "J - K: ", 2) | Out-Null
